# Update cryptos list (Price / Volume(1h) columns, plus a row 34/35 coin swap)
# as produced by the Sun Feb  4 02:44:28 UTC 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: column D ("Price") cells are plain text in this sheet (e.g. "43.030.23"
# uses '.' as a thousands separator). Where the new value still looks like a
# valid number (e.g. "300.80"), force the cell to Text format first so Excel's
# COM layer keeps storing it as a string instead of silently converting it to
# a numeric value.
$ws.Range("D2").Value = "43.030.23"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "2.303.69"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.80"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.22"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.518"
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.35"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.94"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.82"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").Value = "2.662.84"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "2.309.46"
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.783"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").Value = "42.986.67"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.38"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.15"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.19"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.91"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.76"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.03"
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.74"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.78"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Value = "1.999.30"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("E46").Value = "  +1.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.50"
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("E48").Value = "  -3.70%  "
$ws.Range("E49").Value = "  -2.54%  "
$ws.Range("D50").Value = "2.529.96"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.92"
$ws.Range("E51").Value = "  -2.42%  "
